$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.156.78'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '2.844.47'
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '362.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '113.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.572'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +4.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.42%  '

$ws.Range("E13").Value = '  +1.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.79%  '

$ws.Range("D15").Value = '3.287.88'
$ws.Range("E15").Value = '  +2.16%  '

$ws.Range("D16").Value = '2.836.32'
$ws.Range("E16").Value = '  +2.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.919'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.57%  '

$ws.Range("D18").Value = '52.054.73'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000100'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.73%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("E28").Value = '  +2.09%  '

$ws.Range("E30").Value = '  +29.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '53.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.141'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.14%  '

$ws.Range("E35").Value = '  +10.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0846'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.62%  '

$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.59%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.90%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.117'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.06%  '

$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '126.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.98%  '

$ws.Range("E45").Value = '  -3.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.44'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.82%  '

$ws.Range("D47").Value = '2.117.76'
$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("E49").Value = '  +13.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.75%  '
